# aggiornamento 15, 16, 17 marzo
# Append three new daily rows (227-229) to Sheet1, continuing the existing
# date / nuovi pos. / somma mobile 7gg. / somma mobile 7gg. per 100mila
# abitanti series that currently ends at row 226 (A1:D226 -> A1:D229).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A uses a date number format; reuse the format already applied to
# the last existing data row so the new dates render the same way.
$dateFormat = $ws.Cells.Item(226, 1).NumberFormat

$newRows = @(
    @{ Row = 227; Date = 44301; NuoviPos = 1; SommaMobile = 7;  Tasso = 112.1615125781125 },
    @{ Row = 228; Date = 44302; NuoviPos = 4; SommaMobile = 10; Tasso = 160.2307322544464 },
    @{ Row = 229; Date = 44303; NuoviPos = 1; SommaMobile = 8;  Tasso = 128.1845858035571 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 2).Value = $r.NuoviPos
    $ws.Cells.Item($row, 3).Value = $r.SommaMobile
    $ws.Cells.Item($row, 4).Value = $r.Tasso
}
